# Update "Chiffres COVID-19 Valais" daily figures.
# Adds the new day's data (row 206, 2020-09-18) and revises the
# "Nombre de cas / contacts / voyageurs en cours de quarantaine/isolement"
# figures (columns N/O/P) for the preceding days, per the source upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column N (Nombre de cas en cours d'isolement) revisions, rows 181-201 ---
$ws.Range("N181").Value = 69
$ws.Range("N182").Value = 84
$ws.Range("N183").Value = 87
$ws.Range("N184").Value = 87
$ws.Range("N185").Value = 91
$ws.Range("N186").Value = 86
$ws.Range("N187").Value = 84
$ws.Range("N188").Value = 76
$ws.Range("N189").Value = 80
$ws.Range("N190").Value = 79
$ws.Range("N191").Value = 75
$ws.Range("N192").Value = 75

# --- Columns N/O revisions, rows 193-200 ---
$ws.Range("N193").Value = 70
$ws.Range("O193").Value = 229
$ws.Range("N194").Value = 69
$ws.Range("O194").Value = 238
$ws.Range("N195").Value = 82
$ws.Range("O195").Value = 209
$ws.Range("N196").Value = 87
$ws.Range("O196").Value = 210
$ws.Range("N197").Value = 89
$ws.Range("O197").Value = 221
$ws.Range("N198").Value = 94
$ws.Range("O198").Value = 210
$ws.Range("N199").Value = 89
$ws.Range("O199").Value = 219
$ws.Range("N200").Value = 82
$ws.Range("O200").Value = 226

# --- Row 201 ---
$ws.Range("N201").Value = 79

# --- Row 202 (2020-09-14): new-case count + quarantine figures revised ---
$ws.Range("C202").Value = 19
$ws.Range("N202").Value = 90
$ws.Range("P202").Value = 297

# --- Row 203 (2020-09-15) ---
$ws.Range("N203").Value = 90
$ws.Range("O203").Value = 286
$ws.Range("P203").Value = 303

# --- Row 204 (2020-09-16) ---
$ws.Range("C204").Value = 22
$ws.Range("N204").Value = 102
$ws.Range("O204").Value = 327
$ws.Range("P204").Value = 319

# --- Row 205 (2020-09-17) ---
$ws.Range("C205").Value = 12
$ws.Range("N205").Value = 104
$ws.Range("O205").Value = 386
$ws.Range("P205").Value = 302

# --- Row 206 (2020-09-18): brand-new day of data being added ---
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 0
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 8
$ws.Range("I206").Value = 0
$ws.Range("L206").Value = "0"
$ws.Range("M206").Value = "0"
$ws.Range("N206").Value = 102
$ws.Range("O206").Value = 368
$ws.Range("P206").Value = 286

# Restore the view/selection to the top of the frozen pane (B3), matching
# the refreshed sheet view saved with this upload.
$ws.Range("B3").Select()
